# BOM.xlsx update: "Added BBB, sd card & AX12A servos"
#
# - Fixes two shared-string typos (stray leading/trailing spaces)
# - Adds an explicit (formerly implicit/absent) blank A-cell to rows 26-30
# - Adds a fully blank row 31
# - Adds three new BOM rows (33-35): USD-CARDS, BBB (3), AX-12A (1)
# - Nudges the sheet view back to the top-left (A4) from its prior C4 scroll

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---- 1. Fix stray whitespace in two existing shared strings -------------
# C16: "SPOX-3 for AX-12A " -> "SPOX-3 for AX-12A"
$ws.Cells.Item(16, 3).Value = "SPOX-3 for AX-12A"
# B19: " 1.7K RESISTOR" -> "1.7K RESISTOR"
$ws.Cells.Item(19, 2).Value = "1.7K RESISTOR"

# ---- 2. Make column A explicit (blank) on rows 26-30 ---------------------
# These rows already have blank B/C/D cells; add a matching blank A cell so
# every row in the blank-separator/footer block has a full A:D cell set.
foreach ($r in 26..30) {
    $c = $ws.Cells.Item($r, 1)
    $c.Value = "x"
    $c.ClearContents()
}

# ---- 3. Insert a brand-new, fully blank row 31 ---------------------------
foreach ($col in 1..4) {
    $c = $ws.Cells.Item(31, $col)
    $c.Value = "x"
    $c.ClearContents()
}

# ---- 4. Append three new BOM entries (rows 33-35) -------------------------
$ws.Cells.Item(33, 1).Value = "USD-CARDS"
$ws.Cells.Item(33, 4).Value = "http://www.mouser.com/ProductDetail/Apacer/AP16GMCSH4-B/?qs=sGAEpiMZZMtyMAXUUxCBE4AZ7JbBE3hTRlqQ2Hq7Z8o%3d"

$ws.Cells.Item(34, 1).Value = "BBB (3)"
$ws.Cells.Item(34, 4).Value = "http://www.mouser.com/new/embedded-solutions/beagleboneblack/n-5g1kZ2bv0qx"

$ws.Cells.Item(35, 1).Value = "AX-12A (1)"
$ws.Cells.Item(35, 4).Value = "http://www.trossenrobotics.com/dynamixel-ax-12-robot-actuator.aspx"

# ---- 5. Scroll the view back to the top-left (was topLeftCell="C4") ------
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# ---- 6. Best-effort: rename the built-in heading style --------------------
# (Renames "Excel Built-in Heading 1" -> "Excel Built-in Excel Built-in
# Heading 1" to match the authored workbook; harmless no-op if unsupported.)
try {
    $st = $wb.Styles.Item("Excel Built-in Heading 1")
    $st.Name = "Excel Built-in Excel Built-in Heading 1"
} catch {
}
